$d = $word.ActiveDocument

# 1. Wrap the existing "Password" run (first paragraph) in a bookmark,
#    matching <w:bookmarkStart w:id="0" w:name="_Hlk67383494"/> ... <w:bookmarkEnd .../>
$d.Bookmarks.Add("_Hlk67383494", $d.Range(0, 8))

# WordprocessingML namespace used for the raw-XML paragraph inserts below.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 2. New paragraph: "Password" followed by a single space, kept as two
#    separate runs (as in the target markup) by inserting literal OOXML
#    instead of letting adjacent same-format runs coalesce.
$end = $d.Content.End - 1
$d.Range($end, $end).InsertXML("<w:p $wNs><w:r><w:t>Password</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>")

# 3. New paragraph: "Password"
$end = $d.Content.End - 1
$d.Range($end, $end).InsertXML("<w:p $wNs><w:r><w:t>Password</w:t></w:r></w:p>")

# 4. New paragraph: "Password"
$end = $d.Content.End - 1
$d.Range($end, $end).InsertXML("<w:p $wNs><w:r><w:t>Password</w:t></w:r></w:p>")

# 5. New paragraph: "Password"
$end = $d.Content.End - 1
$d.Range($end, $end).InsertXML("<w:p $wNs><w:r><w:t>Password</w:t></w:r></w:p>")
